$p = $ppt.ActivePresentation

# The deck has two duplicate slides (positions 9 and 10) that are exact
# copies of slides 6 and 7. Remove them.
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()
